$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.776.28"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.73%  "
$ws.Range("D3").Value = "'1.869.02"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.00%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "'300.88"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.16%  "
$ws.Range("D6").Value = "'1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.08%  "
$ws.Range("D7").Value = "'0.5324"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Value = "'0.3740"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.00%  "
$ws.Range("D9").Value = "'0.07179"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.72%  "
$ws.Range("D10").Value = "'21.61"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.02%  "
$ws.Range("D11").Value = "'0.8882"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.95%  "
$ws.Range("D12").Value = "'0.08154"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.23%  "
$ws.Range("D13").Value = "'1.876.41"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +14.69%  "
$ws.Range("D14").Value = "'92.84"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Value = "'5.285"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.57%  "
$ws.Range("D16").Value = "'1.001"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.01%  "
$ws.Range("D17").Value = "'14.76"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.10%  "
$ws.Range("D18").Value = "'0.000008475"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.47%  "
$ws.Range("D19").Value = "'1.001"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.08%  "
$ws.Range("D20").Value = "'26.833.41"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.94%  "
$ws.Range("D21").Value = "'4.973"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.02%  "
$ws.Range("E22").Value = "  -1.99%  "
$ws.Range("D23").Value = "'6.379"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.00%  "
$ws.Range("E24").Value = "  -3.23%  "
$ws.Range("D25").Value = "'146.11"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.71%  "
$ws.Range("D26").Value = "'1.733"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.44%  "
$ws.Range("E27").Value = "  -1.30%  "
$ws.Range("D28").Value = "'113.98"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.38%  "
$ws.Range("E29").Value = "  -2.91%  "
$ws.Range("D30").Value = "'4.609"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.29%  "
$ws.Range("D31").Value = "'0.09115"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.50%  "
$ws.Range("D32").Value = "'0.8040"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.21%  "
$ws.Range("D33").Value = "'0.05008"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Value = "'1.165"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.28%  "
$ws.Range("D35").Value = "'2.963"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.38%  "
$ws.Range("D36").Value = "'0.6096"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.60%  "
$ws.Range("D37").Value = "'2.648"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.29%  "
$ws.Range("D38").Value = "'3.198"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.67%  "
$ws.Range("E39").Value = "  -2.65%  "
$ws.Range("D40").Value = "'1.064"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.50%  "
$ws.Range("D41").Value = "'6.500"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.50%  "
$ws.Range("D42").Value = "'8.725"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.51%  "
$ws.Range("D43").Value = "'0.5150"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.73%  "
$ws.Range("D44").Value = "'114.28"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.82%  "
$ws.Range("E45").Value = "  -2.16%  "
$ws.Range("D46").Value = "'1.001"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.06%  "
$ws.Range("D47").Value = "'1.636"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.54%  "
$ws.Range("D48").Value = "'9.932"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.55%  "
$ws.Range("D49").Value = "'37.48"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.28%  "
$ws.Range("D50").Value = "'0.06037"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.28%  "
$ws.Range("D51").Value = "'62.09"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.95%  "
